$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit text
# format first, otherwise Excel auto-converts the assigned string into a
# numeric value (dropping e.g. trailing zeros) instead of keeping literal text.
$ws.Range("D2").Value = "36.623.88"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.069.01"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.56"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.661"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.03"
$ws.Range("E8").Value = "  -7.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.09"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.360"
$ws.Range("E10").Value = "  -7.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0753"
$ws.Range("E11").Value = "  -4.35%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.910"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.69"
$ws.Range("E14").Value = "  -9.23%  "
$ws.Range("D15").Value = "2.367.62"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("D17").Value = "2.085.04"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "36.542.63"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.42"
$ws.Range("E19").Value = "  -13.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.88"
$ws.Range("E20").Value = "  -4.19%  "
$ws.Range("D21").Value = "0.0₃0864"
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.28"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.26"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -4.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.44"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.82"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.61"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.08"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0598"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0818"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("E39").Value = "  -6.42%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.89"
$ws.Range("E40").Value = "  -6.45%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.14"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0216"
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("E44").Value = "  -7.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.55"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").Value = "1.397.63"
$ws.Range("E46").Value = "  +9.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.44"
$ws.Range("E47").Value = "  +8.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.62"
$ws.Range("E48").Value = "  -11.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.36"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.86"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "2.257.12"
$ws.Range("E51").Value = "  +0.48%  "
